$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old remedy rows (rows 3 through 13) entirely, leaving just
# the header (A1) and the single summary row (A2).
$ws.Range("A3:A13").Clear()

# Replace the A2 content with the output filename reference.
$ws.Range("A2").Value = "compliance_remedies_20250521_153619.xlsx"
